$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.803.71'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '3.411.35'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.38%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.411.53'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.569'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.23'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.26%  '
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('E12').Value = '  -2.51%  '
$ws.Range('D13').Value = '3.997.84'
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.09'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('E16').Value = '  -6.80%  '
$ws.Range('D17').Value = '63.867.46'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').Value = '3.395.91'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('E19').Value = '  -2.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.64'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '377.13'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.75'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.96%  '
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('E25').Value = '  -4.81%  '
$ws.Range('E26').Value = '  -0.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.70'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.15%  '
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.06'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.61%  '
$ws.Range('E31').Value = '  -2.93%  '
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '22.87'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.97'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.53'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.57'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E38').Value = '  -1.54%  '
$ws.Range('E39').Value = '  +8.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.16'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.805.02'
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0729'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.83'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.78%  '
$ws.Range('E47').Value = '  -1.83%  '
$ws.Range('E48').Value = '  +11.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '324.03'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.78%  '
$ws.Range('E50').Value = '  -2.42%  '
$ws.Range('E51').Value = '  -2.80%  '
